$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto Price (D) / 1h Volume-change (E) figures - GitHub Actions data pull.
# Source cells are plain text (prices use "." as a thousands separator, e.g. "70.855.96",
# and percents keep significant trailing zeros like "1.00"), so force Text format before
# assigning any value that Excel would otherwise auto-parse as a number.
$ws.Range("D2").Value = "70.855.96"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").Value = "3.586.52"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.88"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "3.579.37"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.621"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.216"
$ws.Range("E10").Value = "  +16.89%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.24"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000320"
$ws.Range("E13").Value = "  +5.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.52"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "4.156.17"
$ws.Range("E15").Value = "  -0.38%  "
$ws.Range("D16").Value = "70.834.36"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.28"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Value = "3.587.32"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.39"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "568.58"
$ws.Range("E20").Value = "  +13.33%  "
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("E22").Value = "  -3.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.57"
$ws.Range("E23").Value = "  -10.02%  "
$ws.Range("E24").Value = "  +5.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.03"
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.82"
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.20"
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.10"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.31"
$ws.Range("E30").Value = "  +1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.22"
$ws.Range("E31").Value = "  -5.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("E32").Value = "  -2.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.114"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "64.14"
$ws.Range("E34").Value = "  -2.97%  "
$ws.Range("E35").Value = "  +3.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "549.43"
$ws.Range("E36").Value = "  -4.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.413"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("D38").Value = "0.0₃0808"
$ws.Range("E38").Value = "  +2.61%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.50"
$ws.Range("E40").Value = "  -2.97%  "
$ws.Range("D41").Value = "3.482.25"
$ws.Range("E41").Value = "  +8.57%  "
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.50"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").Value = "  -2.11%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.02%  "
$ws.Range("E51").Value = "  -4.10%  "
